$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E22").Value = 10
$ws.Range("E24").Value = 5
$ws.Range("F24").Value = " -5 for wrong logic"

$ws.Application.ActiveWindow.ScrollRow = 14
$ws.Range("F21").Select()
